# Rename the "data" worksheet to "Reviews" and update the active selection
# to B18 (matching the author's last cursor position) as per the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")
$ws.Name = "Reviews"

$ws.Activate()
$ws.Range("B18").Select()
